$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.188.98"
$ws.Range("E2").Value = "  +1.53%  "

$ws.Range("D3").Value = "1.781.46"
$ws.Range("E3").Value = "  +0.49%  "

$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.67%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.546"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.29%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.67"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.84%  "

$ws.Range("E9").Value = "  +0.60%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0686"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.47%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0947"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.31%  "

$ws.Range("D12").Value = "2.039.59"
$ws.Range("E12").Value = "  +0.58%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.99"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.82%  "

$ws.Range("D14").Value = "1.781.42"
$ws.Range("E14").Value = "  +0.39%  "

$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "34.111.37"
$ws.Range("E15").Value = "  +1.23%  "

$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.622"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.25%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.13%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.89"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.10%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.25%  "

$ws.Range("D20").Value = "0.0₃0785"
$ws.Range("E20").Value = "  +1.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.53%  "

$ws.Range("E22").Value = "  +0.19%  "

$ws.Range("E23").Value = "  +2.44%  "

$ws.Range("E24").Value = "  -0.51%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.49%  "

$ws.Range("E26").Value = "  +2.46%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.25%  "

$ws.Range("E28").Value = "  +1.63%  "

$ws.Range("E29").Value = "  +0.30%  "

$ws.Range("E30").Value = "  +0.81%  "

$ws.Range("E31").Value = "  +1.67%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.72"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.58%  "

$ws.Range("E33").Value = "  +4.08%  "

$ws.Range("E34").Value = "  -0.20%  "

$ws.Range("D35").Value = "1.441.46"
$ws.Range("E35").Value = "  +4.39%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.653"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.41"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.07%  "

$ws.Range("E38").Value = "  +3.71%  "

$ws.Range("E39").Value = "  +1.14%  "

$ws.Range("E40").Value = "  +0.85%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "80.15"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.79%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.922"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.83%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.67"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.65%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.27%  "

$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0509"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.69%  "

$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.05"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.82%  "

$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0138"
$ws.Range("E47").Value = "  +2.29%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.14%  "

$ws.Range("D49").Value = "1.941.59"
$ws.Range("E49").Value = "  +0.96%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "105.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.82%  "

$ws.Range("E51").Value = "  +0.22%  "
